# Updates crypto price/volume values in the active worksheet to match the
# latest scrape, as produced by the GitHub Actions job on
# Sun Jul 30 05:45:08 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (column D), new Volume(1h) (column E)
$updates = @(
    @{ Row = 2;  D = "29.365.27";      E = "  +0.01%  " },
    @{ Row = 3;  D = "1.880.20";       E = "  +0.33%  " },
    @{ Row = 5;  D = "0.7117";         E = "  -0.05%  " },
    @{ Row = 6;  D = "242.68";         E = "  +0.32%  " },
    @{ Row = 8;  D = "0.08034";        E = "  +3.11%  " },
    @{ Row = 9;  D = "0.3180";         E = "  +2.06%  " },
    @{ Row = 10; D = "25.10";          E = "  -0.41%  " },
    @{ Row = 11; E = "  -1.13%  " },
    @{ Row = 12; D = "1.884.05";       E = "  +0.47%  " },
    @{ Row = 13; D = "5.263";          E = "  +0.46%  " },
    @{ Row = 14; D = "94.97";          E = "  +4.28%  " },
    @{ Row = 15; D = "0.7181";         E = "  +0.65%  " },
    @{ Row = 16; D = "6.360";          E = "  +4.74%  " },
    @{ Row = 17; D = "0.000008619";    E = "  +4.73%  " },
    @{ Row = 18; D = "29.371.93";      E = "  -0.01%  " },
    @{ Row = 19; D = "243.19";         E = "  +0.91%  " },
    @{ Row = 20; D = "2.147.77";       E = "  +0.98%  " },
    @{ Row = 21; D = "13.34";          E = "  +0.65%  " },
    @{ Row = 22; D = "1.002";          E = "  +0.18%  " },
    @{ Row = 23; D = "7.834";          E = "  +0.63%  " },
    @{ Row = 24; E = "  +0.14%  " },
    @{ Row = 25; E = "  -1.24%  " },
    @{ Row = 26; D = "9.106";          E = "  +0.37%  " },
    @{ Row = 27; D = "163.31";         E = "  +0.04%  " },
    @{ Row = 28; E = "  +0.28%  " },
    @{ Row = 29; D = "1.510";          E = "  -0.03%  " },
    @{ Row = 30; D = "4.437";          E = "  +0.37%  " },
    @{ Row = 31; D = "4.347";          E = "  +0.46%  " },
    @{ Row = 32; D = "1.203";          E = "  -6.81%  " },
    @{ Row = 33; E = "  +2.07%  " },
    @{ Row = 34; D = "1.946";          E = "  +0.40%  " },
    @{ Row = 35; D = "0.7735";         E = "  +3.86%  " },
    @{ Row = 36; D = "1.187";          E = "  +0.69%  " },
    @{ Row = 37; D = "2.687";          E = "  -0.44%  " },
    @{ Row = 38; D = "0.01894";        E = "  +1.33%  " },
    @{ Row = 39; D = "1.268.78";       E = "  +3.26%  " },
    @{ Row = 40; D = "2.754";          E = "  +0.90%  " },
    @{ Row = 41; D = "6.502";          E = "  -0.18%  " },
    @{ Row = 42; D = "113.98";         E = "  +2.77%  " },
    @{ Row = 43; D = "0.9095";         E = "  +1.99%  " },
    @{ Row = 44; D = "74.69";          E = "  +2.38%  " },
    @{ Row = 45; E = "  +6.64%  " },
    @{ Row = 46; D = "1.002";          E = "  +0.15%  " },
    @{ Row = 47; D = "2.038.26";       E = "  +0.89%  " },
    @{ Row = 48; E = "  -0.21%  " },
    @{ Row = 49; D = "0.5225" },
    @{ Row = 50; D = "9.542";          E = "  +1.16%  " },
    @{ Row = 51; D = "0.4381";         E = "  +1.34%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($r, 4)
        # Prices are stored as plain text in the source sheet (they use a
        # "1.234.56" grouping style that isn't a real number). Force the
        # cell to text formatting first so Excel doesn't auto-convert the
        # new value into a floating point number and silently round/alter
        # it (e.g. "0.3180" -> 0.318, "25.10" -> 25.1).
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
